$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates.
# Values that look like plain decimal numbers (single '.') must be forced to
# stay text (matching the original inline-string storage) or Excel's COM
# layer will silently coerce them into floating point numbers.
$ws.Range("D2").Value  = "41.637.67"
$ws.Range("D3").Value  = "2.474.13"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value  = "92.40"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value  = "0.552"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0868"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.20"

$ws.Range("D13").Value = "2.855.32"

$ws.Range("D16").Value = "2.474.79"

$ws.Range("D18").Value = "41.588.82"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.73"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.29"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.70"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.83"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.72"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.10"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.29"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"

$ws.Range("D43").Value = "1.987.49"

$ws.Range("D48").Value = "2.712.59"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.03"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.74"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.19"

# Volume(1h) (column E) updates
$ws.Range("E2").Value  = "  +0.16%  "
$ws.Range("E3").Value  = "  +0.50%  "
$ws.Range("E4").Value  = "  -0.03%  "
$ws.Range("E5").Value  = "  +1.32%  "
$ws.Range("E6").Value  = "  +0.81%  "
$ws.Range("E7").Value  = "  +0.63%  "
$ws.Range("E8").Value  = "  +0.03%  "
$ws.Range("E9").Value  = "  +0.66%  "
$ws.Range("E10").Value = "  +9.01%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  +4.37%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +5.50%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("E51").Value = "  +1.05%  "

# Rows 44 and 45 swap places: EnergySwap/VeChain -> VeChain/EnergySwap,
# each with new Price / Volume(1h) values.
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.81"
$ws.Range("E45").Value = "  +1.45%  "
